# Update the cryptos worksheet with refreshed Price / Volume(1h) values.
# (Mirrors the GitHub Actions data-refresh commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" values look like plain numbers (e.g. "606.35"),
# but the source data stores them as text. A leading apostrophe forces
# Excel to keep them as text instead of auto-converting to a number
# (matching how the original cells are stored as inline strings).

$ws.Range("D2").Value = "68.315.11"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.710.12"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'606.35"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "'166.36"
$ws.Range("E6").Value = "  +4.40%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "2.709.13"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "'0.364"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").Value = "'5.30"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'28.47"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "3.204.13"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "68.241.98"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.711.77"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "'11.88"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'370.60"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "'7.62"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Value = "'4.95"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'73.07"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'576.74"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").Value = "'8.17"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'1.60"
$ws.Range("D38").Value = "'162.01"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'19.86"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "'0.377"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'5.39"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D44").Value = "'2.61"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "'0.596"
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("D49").Value = "'154.73"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "'3.91"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +4.49%  "
